$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- Sheet1: termsWithMulitpleLabels ---
$ws1.Range("A279").Value = 'Checked on Oct 21st, 2022, identified following issue:'
$ws1.Range("A279").Font.Bold = $true
$ws1.Range("A279").Font.Color = 0

$ws1.Range("A280").Value = 'sid'
$ws1.Range("B280").Value = 'label_count'
$ws1.Range("C280").Value = 'labels'
$ws1.Range("D280").Value = 'datasets'
$ws1.Range("A281").Value = 'EUPATH_0000001'
$ws1.Range("B281").Value = 2
$ws1.Range("C281").Value = 'Building material | Respiratory and asthma control score (TRACK)'
$ws1.Range("D281").Value = 'icemr_india_cx | gates_crypto | icemr_india_cohort | icemr_india_severe_malaria'
$ws1.Range("A282").Value = 'EUPATH_0010533'
$ws1.Range("B282").Value = 2
$ws1.Range("C282").Value = 'Loose stools in day count | Loose stools in day count, caregiver report'
$ws1.Range("D282").Value = 'gates_crypto | gates_provide | gates_maled'
$ws1.Range("A283").Value = 'EUPATH_0015040'
$ws1.Range("B283").Value = 2
$ws1.Range("C283").Value = 'Unprotected spring | Unprotected dug well or spring'
$ws1.Range("D283").Value = 'gates_gems1a | gates_gems | gates_vida | gates_jilinde_costing_survey'
$ws1.Range("A284").Value = 'EUPATH_0015050'
$ws1.Range("B284").Value = 2
$ws1.Range("C284").Value = 'Bore hole | Protected shallow well or borehole'
$ws1.Range("D284").Value = 'gates_gems1a | gates_gems | gates_vida | gates_jilinde_costing_survey'
$ws1.Range("A285").Value = 'EUPATH_0025047'
$ws1.Range("B285").Value = 2
$ws1.Range("C285").Value = 'Roof material type | Roof material categorization'
$ws1.Range("D285").Value = 'gates_washb_kenya | icemr_malawi | general_promote | gates_shine'
$ws1.Range("A286").Value = 'EUPATH_0036053'
$ws1.Range("B286").Value = 2
$ws1.Range("C286").Value = 'Stool consistencys | Stool consistency, caregiver report'
$ws1.Range("D286").Value = 'gates_crypto | gates_provide'
$ws1.Range("A287").Value = 'EUPATH_0036210'
$ws1.Range("B287").Value = 2
$ws1.Range("C287").Value = 'Household expenditures in local currency | Family expenditures in local currency'
$ws1.Range("D287").Value = 'gates_crypto | gates_provide'
$ws1.Range("A288").Value = 'EUPATH_0051964'
$ws1.Range("B288").Value = 2
$ws1.Range("C288").Value = 'Health care | Healthcare'
$ws1.Range("D288").Value = 'gates_ppfp_choices_kenya_pp | gates_jilinde_costing_survey'
$ws1.Range("A289").Value = 'EUPATH_0054676'
$ws1.Range("B289").Value = 2
$ws1.Range("C289").Value = 'Read newspapers | Read newspaper'
$ws1.Range("D289").Value = 'gates_jilinde_awareness_survey | gates_jilinde_costing_survey | gates_jilinde_demand_creation_evaluation_questionnaire'
$ws1.Range("A290").Value = 'EUPATH_0054747'
$ws1.Range("B290").Value = 2
$ws1.Range("C290").Value = 'Whether listen to radio | Listen to radio'
$ws1.Range("D290").Value = 'gates_jilinde_awareness_survey | gates_jilinde_costing_survey | gates_jilinde_demand_creation_evaluation_questionnaire'
$ws1.Range("A291").Value = 'EUPATH_0054787'
$ws1.Range("B291").Value = 2
$ws1.Range("C291").Value = 'Whether watch TV | Watch television'
$ws1.Range("D291").Value = 'gates_jilinde_awareness_survey | gates_jilinde_costing_survey | gates_jilinde_demand_creation_evaluation_questionnaire'

# --- Sheet2: LabelsUsedMultipleTerms ---
$ws2.Range("A134").Value = 'Checked on Oct 21st, 2022, no issue was found'
$ws2.Range("A134").Font.Bold = $true
$ws2.Range("A134").Font.Color = 0

# --- Sheet3: termWithDifferentParent ---
$ws3.Range("A292").Value = 'Checked on Oct 21st, 2022, identified following issue:'
$ws3.Range("A292").Font.Bold = $true
$ws3.Range("A292").Font.Color = 0

$ws3.Range("A293").Value = 'sid'
$ws3.Range("B293").Value = 'pid_count'
$ws3.Range("C293").Value = 'labels'
$ws3.Range("D293").Value = 'pLabels'
$ws3.Range("E293").Value = 'datasets'
$ws3.Range("F293").Value = "temp"
$ws3.Range("F293").ClearContents()
$ws3.Range("F293").Style = "Normal"

$ws3.Range("A294").Value = 'EUPATH_0015050'
$ws3.Range("B294").Value = 2
$ws3.Range("C294").Value = 'Bore hole|Protected shallow well or borehole'
$ws3.Range("D294").Value = 'Drinking water source|Water source'
$ws3.Range("E294").Value = 'gates_gems1a | gates_gems | gates_vida | gates_jilinde_costing_survey'
$ws3.Range("F294").Value = "temp"
$ws3.Range("F294").ClearContents()
$ws3.Range("F294").Style = "Normal"
$ws3.Range("A295").Value = 'EUPATH_0000001'
$ws3.Range("B295").Value = 2
$ws3.Range("C295").Value = 'Building material|Respiratory and asthma control score (TRACK)'
$ws3.Range("D295").Value = 'Dwelling characteristics|Signs and symptoms'
$ws3.Range("E295").Value = 'icemr_india_cx | gates_crypto | icemr_india_cohort | icemr_india_severe_malaria'
$ws3.Range("F295").Value = "temp"
$ws3.Range("F295").ClearContents()
$ws3.Range("F295").Style = "Normal"
$ws3.Range("A296").Value = 'EUPATH_0049850'
$ws3.Range("B296").Value = 2
$ws3.Range("C296").Value = 'Child vital status'
$ws3.Range("D296").Value = 'Child observation details|Child physical exam'
$ws3.Range("E296").Value = 'gates_ganc | gates_betterbirth'
$ws3.Range("F296").Value = "temp"
$ws3.Range("F296").ClearContents()
$ws3.Range("F296").Style = "Normal"
$ws3.Range("A297").Value = 'EUPATH_0010533'
$ws3.Range("B297").Value = 2
$ws3.Range("C297").Value = 'Loose stools in day count, caregiver report|Loose stools in day count'
$ws3.Range("D297").Value = 'Symptoms|Symptoms by caregiver report'
$ws3.Range("E297").Value = 'gates_crypto | gates_provide | gates_maled'
$ws3.Range("F297").Value = "temp"
$ws3.Range("F297").ClearContents()
$ws3.Range("F297").Style = "Normal"
$ws3.Range("A298").Value = 'EUPATH_0015049'
$ws3.Range("B298").Value = 2
$ws3.Range("C298").Value = 'Other water source'
$ws3.Range("D298").Value = 'Drinking water source|Water source'
$ws3.Range("E298").Value = 'gates_gems1a | gates_gems | gates_vida | gates_jilinde_costing_survey'
$ws3.Range("F298").Value = "temp"
$ws3.Range("F298").ClearContents()
$ws3.Range("F298").Style = "Normal"
$ws3.Range("A299").Value = 'EUPATH_0015038'
$ws3.Range("B299").Value = 2
$ws3.Range("C299").Value = 'Protected spring'
$ws3.Range("D299").Value = 'Drinking water source|Water source'
$ws3.Range("E299").Value = 'gates_gems1a | gates_gems | gates_vida | gates_jilinde_costing_survey'
$ws3.Range("F299").Value = "temp"
$ws3.Range("F299").ClearContents()
$ws3.Range("F299").Style = "Normal"
$ws3.Range("A300").Value = 'EUPATH_0015046'
$ws3.Range("B300").Value = 2
$ws3.Range("C300").Value = 'Rainwater'
$ws3.Range("D300").Value = 'Drinking water source|Water source'
$ws3.Range("E300").Value = 'gates_gems1a | gates_gems | gates_vida | gates_jilinde_costing_survey'
$ws3.Range("F300").Value = "temp"
$ws3.Range("F300").ClearContents()
$ws3.Range("F300").Style = "Normal"
$ws3.Range("A301").Value = 'EUPATH_0054676'
$ws3.Range("B301").Value = 2
$ws3.Range("C301").Value = 'Read newspapers|Read newspaper'
$ws3.Range("D301").Value = 'Activity|Use of media at least once a week'
$ws3.Range("E301").Value = 'gates_jilinde_awareness_survey | gates_jilinde_costing_survey | gates_jilinde_demand_creation_evaluation_questionnaire'
$ws3.Range("F301").Value = "temp"
$ws3.Range("F301").ClearContents()
$ws3.Range("F301").Style = "Normal"
$ws3.Range("A302").Value = 'EUPATH_0036053'
$ws3.Range("B302").Value = 2
$ws3.Range("C302").Value = 'Stool consistency, caregiver report|Stool consistencys'
$ws3.Range("D302").Value = 'Symptoms|Symptoms by caregiver report'
$ws3.Range("E302").Value = 'gates_crypto | gates_provide'
$ws3.Range("F302").Value = "temp"
$ws3.Range("F302").ClearContents()
$ws3.Range("F302").Style = "Normal"
$ws3.Range("A303").Value = 'EUPATH_0031313'
$ws3.Range("B303").Value = 2
$ws3.Range("C303").Value = 'Surface water'
$ws3.Range("D303").Value = 'Drinking water source|Water source'
$ws3.Range("E303").Value = 'gates_vida | gates_jilinde_costing_survey'
$ws3.Range("F303").Value = "temp"
$ws3.Range("F303").ClearContents()
$ws3.Range("F303").Style = "Normal"
$ws3.Range("A304").Value = 'EUPATH_0015040'
$ws3.Range("B304").Value = 2
$ws3.Range("C304").Value = 'Unprotected spring|Unprotected dug well or spring'
$ws3.Range("D304").Value = 'Drinking water source|Water source'
$ws3.Range("E304").Value = 'gates_gems1a | gates_gems | gates_vida | gates_jilinde_costing_survey'
$ws3.Range("F304").Value = "temp"
$ws3.Range("F304").ClearContents()
$ws3.Range("F304").Style = "Normal"
$ws3.Range("A305").Value = 'EUPATH_0054747'
$ws3.Range("B305").Value = 2
$ws3.Range("C305").Value = 'Whether listen to radio|Listen to radio'
$ws3.Range("D305").Value = 'Activity|Use of media at least once a week'
$ws3.Range("E305").Value = 'gates_jilinde_awareness_survey | gates_jilinde_costing_survey | gates_jilinde_demand_creation_evaluation_questionnaire'
$ws3.Range("F305").Value = "temp"
$ws3.Range("F305").ClearContents()
$ws3.Range("F305").Style = "Normal"
$ws3.Range("A306").Value = 'EUPATH_0054787'
$ws3.Range("B306").Value = 2
$ws3.Range("C306").Value = 'Whether watch TV|Watch television'
$ws3.Range("D306").Value = 'Activity|Use of media at least once a week'
$ws3.Range("E306").Value = 'gates_jilinde_awareness_survey | gates_jilinde_costing_survey | gates_jilinde_demand_creation_evaluation_questionnaire'
$ws3.Range("F306").Value = "temp"
$ws3.Range("F306").ClearContents()
$ws3.Range("F306").Style = "Normal"

$ws3.Range("A308").Value = 'Might due to different terms but assigned same ID'

